# Update "想去人数" (F column) and occasional "最低票价" (G column) figures
# on the 展览 sheet and the 全部类型 sheet, reflecting a refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 6870
$ws1.Range("F3").Value = 0
$ws1.Range("F5").Value = 0
$ws1.Range("F6").Value = 151
$ws1.Range("F7").Value = 6587
$ws1.Range("F8").Value = 59
$ws1.Range("F9").Value = 0
$ws1.Range("F10").Value = 1291
$ws1.Range("F14").Value = 0
$ws1.Range("G15").Value = "不可售"
$ws1.Range("F16").Value = 0
$ws1.Range("F19").Value = 4942
$ws1.Range("F20").Value = 95
$ws1.Range("F21").Value = 102
$ws1.Range("F22").Value = 395
$ws1.Range("F23").Value = 0
$ws1.Range("F24").Value = 176

# --- Sheet "全部类型" ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 6870
$ws4.Range("F5").Value = 0
$ws4.Range("F7").Value = 6587
$ws4.Range("F8").Value = 59
$ws4.Range("F12").Value = 0
$ws4.Range("F13").Value = 399
$ws4.Range("F14").Value = 136
$ws4.Range("F15").Value = 18
$ws4.Range("G15").Value = "不可售"
$ws4.Range("F16").Value = 387
$ws4.Range("F18").Value = 9
$ws4.Range("F19").Value = 4942
$ws4.Range("F20").Value = 0
$ws4.Range("F21").Value = 95
$ws4.Range("F22").Value = 102
$ws4.Range("F25").Value = 176
